# Apply the diff: update KNN rows 14-15 with swapped/changed params, add new
# row 16, and move the active selection to G17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Results")

# --- Row 14: 20/[3,5,7] -> 50/5 ---
$ws.Range("D14").Value = 50
$ws.Range("E14").Value = "{'model__n_neighbors': 5}"
$ws.Range("F14").Value = 0.8095238095238095
$ws.Range("G14").Value = 0.720996229973956

# --- Row 15: 50/[3,5,7] -> 20/7 ---
$ws.Range("D15").Value = 20
$ws.Range("E15").Value = "{'model__n_neighbors': 7}"
$ws.Range("F15").Value = 0.8143712574850299
$ws.Range("G15").Value = 0.7369860978207267

# --- Row 16: new KNN / No / UnderSample row ---
$ws.Range("A16").Value = "KNN"
$ws.Range("B16").Value = "No"
$ws.Range("C16").Value = "UnderSample"
$ws.Range("D16").Value = 20
$ws.Range("E16").Value = "{'model__n_neighbors': [3, 5, 7]}"
$ws.Range("F16").Value = 0.7717391304347826
$ws.Range("G16").Value = 0.6199000544466472

# --- Move selection to G17 (matches the updated sheetView) ---
$ws.Range("G17").Select()
